# Remove the erroneous "subgenus" column (AS) from the Materials sheet.
# Darwin Core already carries "genus"; this duplicate/extra "subgenus"
# column (header + "${subgenus}" placeholder row) was flagged during the
# third round of review and is deleted outright here. Excel shifts every
# column to its right (AT.. -> AS..) left by one, which is exactly what
# the native column-delete operation below reproduces, including the
# automatic pruning of the now-unused "subgenus" / "${subgenus}" shared
# strings and the renumbering of shared-string indices referenced by the
# other sheets (e.g. the ExternalLinks sheet's legend row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Column AS (45th column) holds the "subgenus" header in row 1 and the
# "${subgenus}" template value in row 2.
$ws.Columns.Item(45).Delete()
